$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing values ---
$ws.Range("C2").Value = 205
$ws.Range("B6").Value = 2204

# --- New input cell K2 ---
$ws.Range("K2").Value = 10000

# --- Helper formulas in columns H / I for rows 2-4 ---
$ws.Range("H2").Formula = "=B2"
$ws.Range("I2").Formula = "=C2*K2"

$ws.Range("H3").Formula = "=B3"
$ws.Range("I3").Formula = "=C3*K2"

$ws.Range("H4").Formula = "=B6"

# --- Operation labels (row 6 to row 11) ---
$ws.Range("F6").Value = "sub"
$ws.Range("G6").Formula = "=I2"
$ws.Range("H6").Formula = "=I3"
$ws.Range("I6").Formula = "=TRUNC(G6-H6)"

$ws.Range("F7").Value = "sub"
$ws.Range("G7").Formula = "=H2"
$ws.Range("H7").Formula = "=H3"
$ws.Range("I7").Formula = "=TRUNC(G7-H7)"

$ws.Range("F8").Value = "div"
$ws.Range("G8").Formula = "=I6"
$ws.Range("H8").Formula = "=I7"
$ws.Range("I8").Formula = "=TRUNC(G8/H8)"

$ws.Range("F9").Value = "sub"
$ws.Range("G9").Formula = "=H4"
$ws.Range("H9").Formula = "=H3"
$ws.Range("I9").Formula = "=TRUNC(G9-H9)"

$ws.Range("F10").Value = "mul"
$ws.Range("G10").Formula = "=I8"
$ws.Range("H10").Formula = "=I9"
$ws.Range("I10").Formula = "=TRUNC(G10*H10)"

$ws.Range("F11").Value = "add"
$ws.Range("G11").Formula = "=I10"
$ws.Range("H11").Formula = "=I3"
$ws.Range("I11").Formula = "=TRUNC(G11+H11)"

# --- Final result ---
$ws.Range("I13").Formula = "=TRUNC(I11/K2)"

# --- Conditional formatting on G6:I11 ---
$rng = $ws.Range("G6:I11")
$fc = $rng.FormatConditions.Add(1, 5, "65535")
$fc.Font.Color = 0x0006009C
$fc.Interior.Color = 0x00CEC7FF

# --- Selection moved to K2 ---
$ws.Range("K2").Select()
